$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all timestamps in column A (rows 2-97) forward by 2 days
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 2
}

# Update the solar production values for the new date (rows 27-29)
$ws.Cells.Item(27, 2).Value = 3
$ws.Cells.Item(28, 2).Value = 33
$ws.Cells.Item(29, 2).Value = 82
